{"js": "// 1) Remove the existing \"_GoBack\" bookmark. It currently sits right after\n//    \"Laboratorio 12\" at the very start of the document.\nconst doc = context.document;\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Replace \"1234\" with \"5678\" inside the \"Marcar \u201c1234\u201d, ...\" sentence.\nconst body = doc.body;\nconst results = body.search(\"1234\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nconst target = results.items[0];\ntarget.insertText(\"5678\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Re-create the \"_GoBack\" bookmark, collapsed to zero length, right after\n//    the \"5678\" we just inserted (i.e. immediately before the closing\n//    curly quote). `target` now spans exactly the replaced \"5678\" text, so\n//    its end point is exactly where the bookmark needs to go.\nconst insertionPoint = target.getRange(Word.RangeLocation.end);\ninsertionPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the existing \"_GoBack\" bookmark (currently sitting right after\n#    \"Laboratorio 12\" at the top of the document).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Replace \"1234\" with \"5678\" inside the \"Marcar \u201c1234\u201d, ...\" sentence.\n#    Using the positional Execute(...) form so the in/out Range ($rng) is\n#    left spanning exactly the newly-inserted replacement text (\"5678\").\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"1234\"\n$find.MatchCase = $true\n$find.Forward = $true\n$find.Wrap = 0\n$null = $find.Execute(\"1234\", $false, $true, $false, $false, $false, $true, 1, $false, \"5678\", 1)\n\n# 3) Re-create the \"_GoBack\" bookmark, collapsed to zero length, right after\n#    the \"5678\" we just inserted (i.e. immediately before the closing \u201c \u201d).\n#    $rng now covers the replaced \"5678\" text, so its End is exactly there.\n$goBackRange = $d.Range($rng.End, $rng.End)\n$d.Bookmarks.Add(\"_GoBack\", $goBackRange)\n"}
